$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Credentials")

$ws.Range("A11").Value = "jkl"
$ws.Range("B11").Value = "jkl"
$ws.Range("C11").Value = "jkl1"

$ws.Range("A11:C11").Select()
